# Generate Report for Handback
# Adds a new handback row for file "bba96d6b-2231-4154-a1ff-f7bd3ab6182d.md"
# to the Overview sheet, the zh-cn sheet and the de-de sheet, and extends the
# three tables / dimensions accordingly.

$wb = $excel.ActiveWorkbook

$fileId   = "bba96d6b-2231-4154-a1ff-f7bd3ab6182d"
$mdName   = "$fileId.md"
$pathName = "e2e\$mdName"
$statusInSync = "Handed back: in sync with en-US"
$ext = ".md"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("C4").Value = $ext
$wsOverview.Range("E4").Value = $statusInSync
$wsOverview.Range("F4").Value = $statusInSync
$wsOverview.Range("G4").Value = "2016-10-14 07:41:38"
$wsOverview.Range("G4").NumberFormat = $wsOverview.Range("G3").NumberFormat

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f1b2c4d5e6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c/e2e/$mdName", "", "", $pathName) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlf = "$fileId.8222d204b19aad146efd64e63e2a916235e5b757.zh-cn.xlf"

$wsZhCn.Range("B4").Value = $ext
$wsZhCn.Range("C4").Value = $statusInSync
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "True"
$wsZhCn.Range("G4").Value = $zhXlf
$wsZhCn.Range("H4").Value = "2016-10-14 07:41:27"
$wsZhCn.Range("H4").NumberFormat = $wsZhCn.Range("H3").NumberFormat
$wsZhCn.Range("J4").Value = $zhXlf
$wsZhCn.Range("K4").Value = "2016-10-14 07:42:26"
$wsZhCn.Range("K4").NumberFormat = $wsZhCn.Range("K3").NumberFormat
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("O4").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f1b2c4d5e6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c/e2e/$mdName", "", "", $mdName) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c2d3e4f5a6b/e2e/$mdName", "", "", $mdName) | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlf = "$fileId.8222d204b19aad146efd64e63e2a916235e5b757.de-de.xlf"

$wsDeDe.Range("B4").Value = $ext
$wsDeDe.Range("C4").Value = $statusInSync
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "True"
$wsDeDe.Range("G4").Value = $deXlf
$wsDeDe.Range("H4").Value = "2016-10-14 07:41:38"
$wsDeDe.Range("H4").NumberFormat = $wsDeDe.Range("H3").NumberFormat
$wsDeDe.Range("J4").Value = $deXlf
$wsDeDe.Range("K4").Value = "2016-10-14 07:42:43"
$wsDeDe.Range("K4").NumberFormat = $wsDeDe.Range("K3").NumberFormat
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("O4").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f1b2c4d5e6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c/e2e/$mdName", "", "", $mdName) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4b5c6d7e8f9a0b1c2d3e4f5a6b7c8d9e0f1a2b3c/e2e/$mdName", "", "", $mdName) | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P4"))

Write-Output "Handback report row added for $mdName"
